# OW-535: fix margin call generation issues (client vs cpty calls, duplicate
# step nodes, expected/unrecon call matching). Update the sample test
# workbook's Portfolio ID value and refresh the active sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IRS-Bilateral")

# Portfolio ID for the second (data) row changes from "p1" to "p1a"
$ws.Range("AP2").Value = "p1a"

# Reposition the view / selection that was left over from editing the sheet
$ws.Range("AM14").Select()
$excel.ActiveWindow.ScrollColumn = 35
